{"js": "// The source diff adds four package-level parts under customXml/ that\n// SharePoint / OneDrive stamps onto a .docx when it lives in a managed\n// document library with content types + a managed-metadata \"Image Tags\"\n// column:\n//   - customXml/item1.xml       a ct:contentTypeSchema definition\n//   - customXml/item2.xml       a SharePoint FormTemplates declaration\n//   - customXml/itemProps1.xml  the ds:datastoreItem GUID pointer for item1\n//   - customXml/itemProps2.xml  the ds:datastoreItem GUID pointer for item2\n//\n// itemProps1.xml/itemProps2.xml are not independently authored - Word mints\n// one ds:datastoreItem companion part automatically for every custom XML\n// part that is added, so reproducing the edit is two\n// `customXmlParts.add(xml)` calls, one per distinct XML payload.\n\nconst item1Xml = `<?xml version=\"1.0\" encoding=\"utf-8\"?>\n<ct:contentTypeSchema xmlns:ct=\"http://schemas.microsoft.com/office/2006/metadata/contentType\" xmlns:ma=\"http://schemas.microsoft.com/office/2006/metadata/properties/metaAttributes\" ct:_=\"\" ma:_=\"\" ma:contentTypeName=\"Document\" ma:contentTypeID=\"0x01010032561353110FEA42BAAF301666702FCF\" ma:contentTypeVersion=\"13\" ma:contentTypeDescription=\"Create a new document.\" ma:contentTypeScope=\"\" ma:versionID=\"1971c52690794fce746f6f05770e5404\">\n  <xsd:schema xmlns:xsd=\"http://www.w3.org/2001/XMLSchema\" xmlns:xs=\"http://www.w3.org/2001/XMLSchema\" xmlns:p=\"http://schemas.microsoft.com/office/2006/metadata/properties\" xmlns:ns2=\"0aa39948-3de8-40e0-98e7-bf846ca05556\" xmlns:ns3=\"abb1820a-26b6-44f2-b938-9861746a7b9a\" targetNamespace=\"http://schemas.microsoft.com/office/2006/metadata/properties\" ma:root=\"true\" ma:fieldsID=\"2a7e4d74a27da4d52af6b51af46d5402\" ns2:_=\"\" ns3:_=\"\">\n    <xsd:import namespace=\"0aa39948-3de8-40e0-98e7-bf846ca05556\"/>\n    <xsd:import namespace=\"abb1820a-26b6-44f2-b938-9861746a7b9a\"/>\n    <xsd:element name=\"properties\">\n      <xsd:complexType>\n        <xsd:sequence>\n          <xsd:element name=\"documentManagement\">\n            <xsd:complexType>\n              <xsd:all>\n                <xsd:element ref=\"ns2:MediaServiceMetadata\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceFastMetadata\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceAutoKeyPoints\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceKeyPoints\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:lcf76f155ced4ddcb4097134ff3c332f\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns3:TaxCatchAll\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceOCR\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceGenerationTime\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceEventHashCode\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceDateTaken\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns3:SharedWithUsers\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns3:SharedWithDetails\" minOccurs=\"0\"/>\n              </xsd:all>\n            </xsd:complexType>\n          </xsd:element>\n        </xsd:sequence>\n      </xsd:complexType>\n    </xsd:element>\n  </xsd:schema>\n  <xsd:schema xmlns:xsd=\"http://www.w3.org/2001/XMLSchema\" xmlns:xs=\"http://www.w3.org/2001/XMLSchema\" xmlns:dms=\"http://schemas.microsoft.com/office/2006/documentManagement/types\" xmlns:pc=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\" targetNamespace=\"0aa39948-3de8-40e0-98e7-bf846ca05556\" elementFormDefault=\"qualified\">\n    <xsd:import namespace=\"http://schemas.microsoft.com/office/2006/documentManagement/types\"/>\n    <xsd:import namespace=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\"/>\n    <xsd:element name=\"MediaServiceMetadata\" ma:index=\"8\" nillable=\"true\" ma:displayName=\"MediaServiceMetadata\" ma:hidden=\"true\" ma:internalName=\"MediaServiceMetadata\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Note\"/>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceFastMetadata\" ma:index=\"9\" nillable=\"true\" ma:displayName=\"MediaServiceFastMetadata\" ma:hidden=\"true\" ma:internalName=\"MediaServiceFastMetadata\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Note\"/>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceAutoKeyPoints\" ma:index=\"10\" nillable=\"true\" ma:displayName=\"MediaServiceAutoKeyPoints\" ma:hidden=\"true\" ma:internalName=\"MediaServiceAutoKeyPoints\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Note\"/>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceKeyPoints\" ma:index=\"11\" nillable=\"true\" ma:displayName=\"KeyPoints\" ma:internalName=\"MediaServiceKeyPoints\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Note\">\n          <xsd:maxLength value=\"255\"/>\n        </xsd:restriction>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"lcf76f155ced4ddcb4097134ff3c332f\" ma:index=\"13\" nillable=\"true\" ma:taxonomy=\"true\" ma:internalName=\"lcf76f155ced4ddcb4097134ff3c332f\" ma:taxonomyFieldName=\"MediaServiceImageTags\" ma:displayName=\"Image Tags\" ma:readOnly=\"false\" ma:fieldId=\"{5cf76f15-5ced-4ddc-b409-7134ff3c332f}\" ma:taxonomyMulti=\"true\" ma:sspId=\"dd802298-ac7f-4dc9-a73d-133dd7ac0fd3\" ma:termSetId=\"09814cd3-568e-fe90-9814-8d621ff8fb84\" ma:anchorId=\"fba54fb3-c3e1-fe81-a776-ca4b69148c4d\" ma:open=\"true\" ma:isKeyword=\"false\">\n      <xsd:complexType>\n        <xsd:sequence>\n          <xsd:element ref=\"pc:Terms\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        </xsd:sequence>\n      </xsd:complexType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceOCR\" ma:index=\"15\" nillable=\"true\" ma:displayName=\"Extracted Text\" ma:internalName=\"MediaServiceOCR\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Note\">\n          <xsd:maxLength value=\"255\"/>\n        </xsd:restriction>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceGenerationTime\" ma:index=\"16\" nillable=\"true\" ma:displayName=\"MediaServiceGenerationTime\" ma:hidden=\"true\" ma:internalName=\"MediaServiceGenerationTime\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Text\"/>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceEventHashCode\" ma:index=\"17\" nillable=\"true\" ma:displayName=\"MediaServiceEventHashCode\" ma:hidden=\"true\" ma:internalName=\"MediaServiceEventHashCode\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Text\"/>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceDateTaken\" ma:index=\"18\" nillable=\"true\" ma:displayName=\"MediaServiceDateTaken\" ma:hidden=\"true\" ma:internalName=\"MediaServiceDateTaken\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Text\"/>\n      </xsd:simpleType>\n    </xsd:element>\n  </xsd:schema>\n  <xsd:schema xmlns:xsd=\"http://www.w3.org/2001/XMLSchema\" xmlns:xs=\"http://www.w3.org/2001/XMLSchema\" xmlns:dms=\"http://schemas.microsoft.com/office/2006/documentManagement/types\" xmlns:pc=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\" targetNamespace=\"abb1820a-26b6-44f2-b938-9861746a7b9a\" elementFormDefault=\"qualified\">\n    <xsd:import namespace=\"http://schemas.microsoft.com/office/2006/documentManagement/types\"/>\n    <xsd:import namespace=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\"/>\n    <xsd:element name=\"TaxCatchAll\" ma:index=\"14\" nillable=\"true\" ma:displayName=\"Taxonomy Catch All Column\" ma:hidden=\"true\" ma:list=\"{dd007225-e316-4aa6-b526-431e5738c9a1}\" ma:internalName=\"TaxCatchAll\" ma:showField=\"CatchAllData\" ma:web=\"abb1820a-26b6-44f2-b938-9861746a7b9a\">\n      <xsd:complexType>\n        <xsd:complexContent>\n          <xsd:extension base=\"dms:MultiChoiceLookup\">\n            <xsd:sequence>\n              <xsd:element name=\"Value\" type=\"dms:Lookup\" maxOccurs=\"unbounded\" minOccurs=\"0\" nillable=\"true\"/>\n            </xsd:sequence>\n          </xsd:extension>\n        </xsd:complexContent>\n      </xsd:complexType>\n    </xsd:element>\n    <xsd:element name=\"SharedWithUsers\" ma:index=\"19\" nillable=\"true\" ma:displayName=\"Shared With\" ma:internalName=\"SharedWithUsers\" ma:readOnly=\"true\">\n      <xsd:complexType>\n        <xsd:complexContent>\n          <xsd:extension base=\"dms:UserMulti\">\n            <xsd:sequence>\n              <xsd:element name=\"UserInfo\" minOccurs=\"0\" maxOccurs=\"unbounded\">\n                <xsd:complexType>\n                  <xsd:sequence>\n                    <xsd:element name=\"DisplayName\" type=\"xsd:string\" minOccurs=\"0\"/>\n                    <xsd:element name=\"AccountId\" type=\"dms:UserId\" minOccurs=\"0\" nillable=\"true\"/>\n                    <xsd:element name=\"AccountType\" type=\"xsd:string\" minOccurs=\"0\"/>\n                  </xsd:sequence>\n                </xsd:complexType>\n              </xsd:element>\n            </xsd:sequence>\n          </xsd:extension>\n        </xsd:complexContent>\n      </xsd:complexType>\n    </xsd:element>\n    <xsd:element name=\"SharedWithDetails\" ma:index=\"20\" nillable=\"true\" ma:displayName=\"Shared With Details\" ma:internalName=\"SharedWithDetails\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Note\">\n          <xsd:maxLength value=\"255\"/>\n        </xsd:restriction>\n      </xsd:simpleType>\n    </xsd:element>\n  </xsd:schema>\n  <xsd:schema xmlns=\"http://schemas.openxmlformats.org/package/2006/metadata/core-properties\" xmlns:xsd=\"http://www.w3.org/2001/XMLSchema\" xmlns:xsi=\"http://www.w3.org/2001/XMLSchema-instance\" xmlns:dc=\"http://purl.org/dc/elements/1.1/\" xmlns:dcterms=\"http://purl.org/dc/terms/\" xmlns:odoc=\"http://schemas.microsoft.com/internal/obd\" targetNamespace=\"http://schemas.openxmlformats.org/package/2006/metadata/core-properties\" elementFormDefault=\"qualified\" attributeFormDefault=\"unqualified\" blockDefault=\"#all\">\n    <xsd:import namespace=\"http://purl.org/dc/elements/1.1/\" schemaLocation=\"http://dublincore.org/schemas/xmls/qdc/2003/04/02/dc.xsd\"/>\n    <xsd:import namespace=\"http://purl.org/dc/terms/\" schemaLocation=\"http://dublincore.org/schemas/xmls/qdc/2003/04/02/dcterms.xsd\"/>\n    <xsd:element name=\"coreProperties\" type=\"CT_coreProperties\"/>\n    <xsd:complexType name=\"CT_coreProperties\">\n      <xsd:all>\n        <xsd:element ref=\"dc:creator\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element ref=\"dcterms:created\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element ref=\"dc:identifier\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element name=\"contentType\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\" ma:index=\"0\" ma:displayName=\"Content Type\"/>\n        <xsd:element ref=\"dc:title\" minOccurs=\"0\" maxOccurs=\"1\" ma:index=\"4\" ma:displayName=\"Title\"/>\n        <xsd:element ref=\"dc:subject\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element ref=\"dc:description\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element name=\"keywords\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\"/>\n        <xsd:element ref=\"dc:language\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element name=\"category\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\"/>\n        <xsd:element name=\"version\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\"/>\n        <xsd:element name=\"revision\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\">\n          <xsd:annotation>\n            <xsd:documentation>\n                        This value indicates the number of saves or revisions. The application is responsible for updating this value after each revision.\n                    </xsd:documentation>\n          </xsd:annotation>\n        </xsd:element>\n        <xsd:element name=\"lastModifiedBy\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\"/>\n        <xsd:element ref=\"dcterms:modified\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element name=\"contentStatus\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\"/>\n      </xsd:all>\n    </xsd:complexType>\n  </xsd:schema>\n  <xs:schema xmlns:pc=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\" xmlns:xs=\"http://www.w3.org/2001/XMLSchema\" targetNamespace=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\" elementFormDefault=\"qualified\" attributeFormDefault=\"unqualified\">\n    <xs:element name=\"Person\">\n      <xs:complexType>\n        <xs:sequence>\n          <xs:element ref=\"pc:DisplayName\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:AccountId\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:AccountType\" minOccurs=\"0\"/>\n        </xs:sequence>\n      </xs:complexType>\n    </xs:element>\n    <xs:element name=\"DisplayName\" type=\"xs:string\"/>\n    <xs:element name=\"AccountId\" type=\"xs:string\"/>\n    <xs:element name=\"AccountType\" type=\"xs:string\"/>\n    <xs:element name=\"BDCAssociatedEntity\">\n      <xs:complexType>\n        <xs:sequence>\n          <xs:element ref=\"pc:BDCEntity\" minOccurs=\"0\" maxOccurs=\"unbounded\"/>\n        </xs:sequence>\n        <xs:attribute ref=\"pc:EntityNamespace\"/>\n        <xs:attribute ref=\"pc:EntityName\"/>\n        <xs:attribute ref=\"pc:SystemInstanceName\"/>\n        <xs:attribute ref=\"pc:AssociationName\"/>\n      </xs:complexType>\n    </xs:element>\n    <xs:attribute name=\"EntityNamespace\" type=\"xs:string\"/>\n    <xs:attribute name=\"EntityName\" type=\"xs:string\"/>\n    <xs:attribute name=\"SystemInstanceName\" type=\"xs:string\"/>\n    <xs:attribute name=\"AssociationName\" type=\"xs:string\"/>\n    <xs:element name=\"BDCEntity\">\n      <xs:complexType>\n        <xs:sequence>\n          <xs:element ref=\"pc:EntityDisplayName\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:EntityInstanceReference\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:EntityId1\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:EntityId2\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:EntityId3\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:EntityId4\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:EntityId5\" minOccurs=\"0\"/>\n        </xs:sequence>\n      </xs:complexType>\n    </xs:element>\n    <xs:element name=\"EntityDisplayName\" type=\"xs:string\"/>\n    <xs:element name=\"EntityInstanceReference\" type=\"xs:string\"/>\n    <xs:element name=\"EntityId1\" type=\"xs:string\"/>\n    <xs:element name=\"EntityId2\" type=\"xs:string\"/>\n    <xs:element name=\"EntityId3\" type=\"xs:string\"/>\n    <xs:element name=\"EntityId4\" type=\"xs:string\"/>\n    <xs:element name=\"EntityId5\" type=\"xs:string\"/>\n    <xs:element name=\"Terms\">\n      <xs:complexType>\n        <xs:sequence>\n          <xs:element ref=\"pc:TermInfo\" minOccurs=\"0\" maxOccurs=\"unbounded\"/>\n        </xs:sequence>\n      </xs:complexType>\n    </xs:element>\n    <xs:element name=\"TermInfo\">\n      <xs:complexType>\n        <xs:sequence>\n          <xs:element ref=\"pc:TermName\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:TermId\" minOccurs=\"0\"/>\n        </xs:sequence>\n      </xs:complexType>\n    </xs:element>\n    <xs:element name=\"TermName\" type=\"xs:string\"/>\n    <xs:element name=\"TermId\" type=\"xs:string\"/>\n  </xs:schema>\n</ct:contentTypeSchema>`;\n\nconst item2Xml = `<?xml version=\"1.0\" encoding=\"utf-8\"?>\n<?mso-contentType ?>\n<FormTemplates xmlns=\"http://schemas.microsoft.com/sharepoint/v3/contenttype/forms\">\n  <Display>DocumentLibraryForm</Display>\n  <Edit>DocumentLibraryForm</Edit>\n  <New>DocumentLibraryForm</New>\n</FormTemplates>`;\n\nconst customXmlParts = context.document.customXmlParts;\ncustomXmlParts.add(item1Xml);\ncustomXmlParts.add(item2Xml);\n\nawait context.sync();\n", "ps1": "# The source diff adds four package-level parts under customXml/ that\n# SharePoint / OneDrive stamps onto a .docx when it lives in a managed\n# document library with content types + a managed-metadata \"Image Tags\"\n# column:\n#   - customXml/item1.xml       a ct:contentTypeSchema definition\n#   - customXml/item2.xml       a SharePoint FormTemplates declaration\n#   - customXml/itemProps1.xml  the ds:datastoreItem GUID pointer for item1\n#   - customXml/itemProps2.xml  the ds:datastoreItem GUID pointer for item2\n#\n# itemProps1.xml/itemProps2.xml are not independently authored - Word mints\n# one ds:datastoreItem companion part automatically for every custom XML\n# part that is added, so reproducing the edit is two\n# CustomXMLParts.Add(xml) calls, one per distinct XML payload.\n\n$d = $word.ActiveDocument\n\n$item1Xml = @'\n<?xml version=\"1.0\" encoding=\"utf-8\"?>\n<ct:contentTypeSchema xmlns:ct=\"http://schemas.microsoft.com/office/2006/metadata/contentType\" xmlns:ma=\"http://schemas.microsoft.com/office/2006/metadata/properties/metaAttributes\" ct:_=\"\" ma:_=\"\" ma:contentTypeName=\"Document\" ma:contentTypeID=\"0x01010032561353110FEA42BAAF301666702FCF\" ma:contentTypeVersion=\"13\" ma:contentTypeDescription=\"Create a new document.\" ma:contentTypeScope=\"\" ma:versionID=\"1971c52690794fce746f6f05770e5404\">\n  <xsd:schema xmlns:xsd=\"http://www.w3.org/2001/XMLSchema\" xmlns:xs=\"http://www.w3.org/2001/XMLSchema\" xmlns:p=\"http://schemas.microsoft.com/office/2006/metadata/properties\" xmlns:ns2=\"0aa39948-3de8-40e0-98e7-bf846ca05556\" xmlns:ns3=\"abb1820a-26b6-44f2-b938-9861746a7b9a\" targetNamespace=\"http://schemas.microsoft.com/office/2006/metadata/properties\" ma:root=\"true\" ma:fieldsID=\"2a7e4d74a27da4d52af6b51af46d5402\" ns2:_=\"\" ns3:_=\"\">\n    <xsd:import namespace=\"0aa39948-3de8-40e0-98e7-bf846ca05556\"/>\n    <xsd:import namespace=\"abb1820a-26b6-44f2-b938-9861746a7b9a\"/>\n    <xsd:element name=\"properties\">\n      <xsd:complexType>\n        <xsd:sequence>\n          <xsd:element name=\"documentManagement\">\n            <xsd:complexType>\n              <xsd:all>\n                <xsd:element ref=\"ns2:MediaServiceMetadata\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceFastMetadata\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceAutoKeyPoints\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceKeyPoints\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:lcf76f155ced4ddcb4097134ff3c332f\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns3:TaxCatchAll\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceOCR\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceGenerationTime\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceEventHashCode\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns2:MediaServiceDateTaken\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns3:SharedWithUsers\" minOccurs=\"0\"/>\n                <xsd:element ref=\"ns3:SharedWithDetails\" minOccurs=\"0\"/>\n              </xsd:all>\n            </xsd:complexType>\n          </xsd:element>\n        </xsd:sequence>\n      </xsd:complexType>\n    </xsd:element>\n  </xsd:schema>\n  <xsd:schema xmlns:xsd=\"http://www.w3.org/2001/XMLSchema\" xmlns:xs=\"http://www.w3.org/2001/XMLSchema\" xmlns:dms=\"http://schemas.microsoft.com/office/2006/documentManagement/types\" xmlns:pc=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\" targetNamespace=\"0aa39948-3de8-40e0-98e7-bf846ca05556\" elementFormDefault=\"qualified\">\n    <xsd:import namespace=\"http://schemas.microsoft.com/office/2006/documentManagement/types\"/>\n    <xsd:import namespace=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\"/>\n    <xsd:element name=\"MediaServiceMetadata\" ma:index=\"8\" nillable=\"true\" ma:displayName=\"MediaServiceMetadata\" ma:hidden=\"true\" ma:internalName=\"MediaServiceMetadata\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Note\"/>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceFastMetadata\" ma:index=\"9\" nillable=\"true\" ma:displayName=\"MediaServiceFastMetadata\" ma:hidden=\"true\" ma:internalName=\"MediaServiceFastMetadata\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Note\"/>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceAutoKeyPoints\" ma:index=\"10\" nillable=\"true\" ma:displayName=\"MediaServiceAutoKeyPoints\" ma:hidden=\"true\" ma:internalName=\"MediaServiceAutoKeyPoints\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Note\"/>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceKeyPoints\" ma:index=\"11\" nillable=\"true\" ma:displayName=\"KeyPoints\" ma:internalName=\"MediaServiceKeyPoints\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Note\">\n          <xsd:maxLength value=\"255\"/>\n        </xsd:restriction>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"lcf76f155ced4ddcb4097134ff3c332f\" ma:index=\"13\" nillable=\"true\" ma:taxonomy=\"true\" ma:internalName=\"lcf76f155ced4ddcb4097134ff3c332f\" ma:taxonomyFieldName=\"MediaServiceImageTags\" ma:displayName=\"Image Tags\" ma:readOnly=\"false\" ma:fieldId=\"{5cf76f15-5ced-4ddc-b409-7134ff3c332f}\" ma:taxonomyMulti=\"true\" ma:sspId=\"dd802298-ac7f-4dc9-a73d-133dd7ac0fd3\" ma:termSetId=\"09814cd3-568e-fe90-9814-8d621ff8fb84\" ma:anchorId=\"fba54fb3-c3e1-fe81-a776-ca4b69148c4d\" ma:open=\"true\" ma:isKeyword=\"false\">\n      <xsd:complexType>\n        <xsd:sequence>\n          <xsd:element ref=\"pc:Terms\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        </xsd:sequence>\n      </xsd:complexType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceOCR\" ma:index=\"15\" nillable=\"true\" ma:displayName=\"Extracted Text\" ma:internalName=\"MediaServiceOCR\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Note\">\n          <xsd:maxLength value=\"255\"/>\n        </xsd:restriction>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceGenerationTime\" ma:index=\"16\" nillable=\"true\" ma:displayName=\"MediaServiceGenerationTime\" ma:hidden=\"true\" ma:internalName=\"MediaServiceGenerationTime\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Text\"/>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceEventHashCode\" ma:index=\"17\" nillable=\"true\" ma:displayName=\"MediaServiceEventHashCode\" ma:hidden=\"true\" ma:internalName=\"MediaServiceEventHashCode\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Text\"/>\n      </xsd:simpleType>\n    </xsd:element>\n    <xsd:element name=\"MediaServiceDateTaken\" ma:index=\"18\" nillable=\"true\" ma:displayName=\"MediaServiceDateTaken\" ma:hidden=\"true\" ma:internalName=\"MediaServiceDateTaken\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Text\"/>\n      </xsd:simpleType>\n    </xsd:element>\n  </xsd:schema>\n  <xsd:schema xmlns:xsd=\"http://www.w3.org/2001/XMLSchema\" xmlns:xs=\"http://www.w3.org/2001/XMLSchema\" xmlns:dms=\"http://schemas.microsoft.com/office/2006/documentManagement/types\" xmlns:pc=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\" targetNamespace=\"abb1820a-26b6-44f2-b938-9861746a7b9a\" elementFormDefault=\"qualified\">\n    <xsd:import namespace=\"http://schemas.microsoft.com/office/2006/documentManagement/types\"/>\n    <xsd:import namespace=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\"/>\n    <xsd:element name=\"TaxCatchAll\" ma:index=\"14\" nillable=\"true\" ma:displayName=\"Taxonomy Catch All Column\" ma:hidden=\"true\" ma:list=\"{dd007225-e316-4aa6-b526-431e5738c9a1}\" ma:internalName=\"TaxCatchAll\" ma:showField=\"CatchAllData\" ma:web=\"abb1820a-26b6-44f2-b938-9861746a7b9a\">\n      <xsd:complexType>\n        <xsd:complexContent>\n          <xsd:extension base=\"dms:MultiChoiceLookup\">\n            <xsd:sequence>\n              <xsd:element name=\"Value\" type=\"dms:Lookup\" maxOccurs=\"unbounded\" minOccurs=\"0\" nillable=\"true\"/>\n            </xsd:sequence>\n          </xsd:extension>\n        </xsd:complexContent>\n      </xsd:complexType>\n    </xsd:element>\n    <xsd:element name=\"SharedWithUsers\" ma:index=\"19\" nillable=\"true\" ma:displayName=\"Shared With\" ma:internalName=\"SharedWithUsers\" ma:readOnly=\"true\">\n      <xsd:complexType>\n        <xsd:complexContent>\n          <xsd:extension base=\"dms:UserMulti\">\n            <xsd:sequence>\n              <xsd:element name=\"UserInfo\" minOccurs=\"0\" maxOccurs=\"unbounded\">\n                <xsd:complexType>\n                  <xsd:sequence>\n                    <xsd:element name=\"DisplayName\" type=\"xsd:string\" minOccurs=\"0\"/>\n                    <xsd:element name=\"AccountId\" type=\"dms:UserId\" minOccurs=\"0\" nillable=\"true\"/>\n                    <xsd:element name=\"AccountType\" type=\"xsd:string\" minOccurs=\"0\"/>\n                  </xsd:sequence>\n                </xsd:complexType>\n              </xsd:element>\n            </xsd:sequence>\n          </xsd:extension>\n        </xsd:complexContent>\n      </xsd:complexType>\n    </xsd:element>\n    <xsd:element name=\"SharedWithDetails\" ma:index=\"20\" nillable=\"true\" ma:displayName=\"Shared With Details\" ma:internalName=\"SharedWithDetails\" ma:readOnly=\"true\">\n      <xsd:simpleType>\n        <xsd:restriction base=\"dms:Note\">\n          <xsd:maxLength value=\"255\"/>\n        </xsd:restriction>\n      </xsd:simpleType>\n    </xsd:element>\n  </xsd:schema>\n  <xsd:schema xmlns=\"http://schemas.openxmlformats.org/package/2006/metadata/core-properties\" xmlns:xsd=\"http://www.w3.org/2001/XMLSchema\" xmlns:xsi=\"http://www.w3.org/2001/XMLSchema-instance\" xmlns:dc=\"http://purl.org/dc/elements/1.1/\" xmlns:dcterms=\"http://purl.org/dc/terms/\" xmlns:odoc=\"http://schemas.microsoft.com/internal/obd\" targetNamespace=\"http://schemas.openxmlformats.org/package/2006/metadata/core-properties\" elementFormDefault=\"qualified\" attributeFormDefault=\"unqualified\" blockDefault=\"#all\">\n    <xsd:import namespace=\"http://purl.org/dc/elements/1.1/\" schemaLocation=\"http://dublincore.org/schemas/xmls/qdc/2003/04/02/dc.xsd\"/>\n    <xsd:import namespace=\"http://purl.org/dc/terms/\" schemaLocation=\"http://dublincore.org/schemas/xmls/qdc/2003/04/02/dcterms.xsd\"/>\n    <xsd:element name=\"coreProperties\" type=\"CT_coreProperties\"/>\n    <xsd:complexType name=\"CT_coreProperties\">\n      <xsd:all>\n        <xsd:element ref=\"dc:creator\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element ref=\"dcterms:created\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element ref=\"dc:identifier\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element name=\"contentType\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\" ma:index=\"0\" ma:displayName=\"Content Type\"/>\n        <xsd:element ref=\"dc:title\" minOccurs=\"0\" maxOccurs=\"1\" ma:index=\"4\" ma:displayName=\"Title\"/>\n        <xsd:element ref=\"dc:subject\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element ref=\"dc:description\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element name=\"keywords\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\"/>\n        <xsd:element ref=\"dc:language\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element name=\"category\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\"/>\n        <xsd:element name=\"version\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\"/>\n        <xsd:element name=\"revision\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\">\n          <xsd:annotation>\n            <xsd:documentation>\n                        This value indicates the number of saves or revisions. The application is responsible for updating this value after each revision.\n                    </xsd:documentation>\n          </xsd:annotation>\n        </xsd:element>\n        <xsd:element name=\"lastModifiedBy\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\"/>\n        <xsd:element ref=\"dcterms:modified\" minOccurs=\"0\" maxOccurs=\"1\"/>\n        <xsd:element name=\"contentStatus\" minOccurs=\"0\" maxOccurs=\"1\" type=\"xsd:string\"/>\n      </xsd:all>\n    </xsd:complexType>\n  </xsd:schema>\n  <xs:schema xmlns:pc=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\" xmlns:xs=\"http://www.w3.org/2001/XMLSchema\" targetNamespace=\"http://schemas.microsoft.com/office/infopath/2007/PartnerControls\" elementFormDefault=\"qualified\" attributeFormDefault=\"unqualified\">\n    <xs:element name=\"Person\">\n      <xs:complexType>\n        <xs:sequence>\n          <xs:element ref=\"pc:DisplayName\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:AccountId\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:AccountType\" minOccurs=\"0\"/>\n        </xs:sequence>\n      </xs:complexType>\n    </xs:element>\n    <xs:element name=\"DisplayName\" type=\"xs:string\"/>\n    <xs:element name=\"AccountId\" type=\"xs:string\"/>\n    <xs:element name=\"AccountType\" type=\"xs:string\"/>\n    <xs:element name=\"BDCAssociatedEntity\">\n      <xs:complexType>\n        <xs:sequence>\n          <xs:element ref=\"pc:BDCEntity\" minOccurs=\"0\" maxOccurs=\"unbounded\"/>\n        </xs:sequence>\n        <xs:attribute ref=\"pc:EntityNamespace\"/>\n        <xs:attribute ref=\"pc:EntityName\"/>\n        <xs:attribute ref=\"pc:SystemInstanceName\"/>\n        <xs:attribute ref=\"pc:AssociationName\"/>\n      </xs:complexType>\n    </xs:element>\n    <xs:attribute name=\"EntityNamespace\" type=\"xs:string\"/>\n    <xs:attribute name=\"EntityName\" type=\"xs:string\"/>\n    <xs:attribute name=\"SystemInstanceName\" type=\"xs:string\"/>\n    <xs:attribute name=\"AssociationName\" type=\"xs:string\"/>\n    <xs:element name=\"BDCEntity\">\n      <xs:complexType>\n        <xs:sequence>\n          <xs:element ref=\"pc:EntityDisplayName\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:EntityInstanceReference\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:EntityId1\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:EntityId2\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:EntityId3\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:EntityId4\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:EntityId5\" minOccurs=\"0\"/>\n        </xs:sequence>\n      </xs:complexType>\n    </xs:element>\n    <xs:element name=\"EntityDisplayName\" type=\"xs:string\"/>\n    <xs:element name=\"EntityInstanceReference\" type=\"xs:string\"/>\n    <xs:element name=\"EntityId1\" type=\"xs:string\"/>\n    <xs:element name=\"EntityId2\" type=\"xs:string\"/>\n    <xs:element name=\"EntityId3\" type=\"xs:string\"/>\n    <xs:element name=\"EntityId4\" type=\"xs:string\"/>\n    <xs:element name=\"EntityId5\" type=\"xs:string\"/>\n    <xs:element name=\"Terms\">\n      <xs:complexType>\n        <xs:sequence>\n          <xs:element ref=\"pc:TermInfo\" minOccurs=\"0\" maxOccurs=\"unbounded\"/>\n        </xs:sequence>\n      </xs:complexType>\n    </xs:element>\n    <xs:element name=\"TermInfo\">\n      <xs:complexType>\n        <xs:sequence>\n          <xs:element ref=\"pc:TermName\" minOccurs=\"0\"/>\n          <xs:element ref=\"pc:TermId\" minOccurs=\"0\"/>\n        </xs:sequence>\n      </xs:complexType>\n    </xs:element>\n    <xs:element name=\"TermName\" type=\"xs:string\"/>\n    <xs:element name=\"TermId\" type=\"xs:string\"/>\n  </xs:schema>\n</ct:contentTypeSchema>\n'@\n\n$item2Xml = @'\n<?xml version=\"1.0\" encoding=\"utf-8\"?>\n<?mso-contentType ?>\n<FormTemplates xmlns=\"http://schemas.microsoft.com/sharepoint/v3/contenttype/forms\">\n  <Display>DocumentLibraryForm</Display>\n  <Edit>DocumentLibraryForm</Edit>\n  <New>DocumentLibraryForm</New>\n</FormTemplates>\n'@\n\n$d.CustomXMLParts.Add($item1Xml) | Out-Null\n$d.CustomXMLParts.Add($item2Xml) | Out-Null\n"}
